$d = $word.ActiveDocument

$replacements = @(
    @("990÷8=", "472÷9="),
    @("214÷3=", "943÷7="),
    @("368÷9=", "250÷7="),
    @("512÷2=", "555÷6="),
    @("708÷2=", "431÷5="),
    @("427÷4=", "282÷8="),
    @("397÷4=", "600÷7="),
    @("840÷3=", "417÷8="),
    @("511÷7=", "787÷5="),
    @("240÷6=", "637÷2="),
    @("372÷7=", "949÷4="),
    @("548÷8=", "697÷9="),
    @("538÷7=", "642÷6="),
    @("810÷8=", "668÷8="),
    @("224÷5=", "519÷9="),
    @("366÷9=", "451÷8="),
    @("236÷9=", "983÷6="),
    @("847÷2=", "577÷8="),
    @("194÷2=", "319÷6="),
    @("439÷4=", "388÷9="),
    @("776÷7=", "916÷7="),
    @("562÷4=", "196÷3="),
    @("318÷7=", "755÷4="),
    @("186÷4=", "493÷2="),
    @("763÷5=", "405÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
